$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-7: price/volume updates (no reorder)
$ws.Range("D2").Value = "'305.68"
$ws.Range("E2").Value = "'1.77%"
$ws.Range("D3").Value = "'36.22"
$ws.Range("E3").Value = "'-5.03%"
$ws.Range("D4").Value = "'5.033"
$ws.Range("E4").Value = "'0.38%"
$ws.Range("D5").Value = "'0.07828"
$ws.Range("E5").Value = "'1.43%"
$ws.Range("D6").Value = "'2.120"
$ws.Range("E6").Value = "'-3.00%"
$ws.Range("D7").Value = "'7.922"
$ws.Range("E7").Value = "'-0.52%"

# Rows 8-17: coin list shifted down by one (GateToken moved to top) + price/volume updates
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.066"
$ws.Range("E8").Value = "'1.79%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9183"
$ws.Range("E9").Value = "'0.10%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.09567"
$ws.Range("E10").Value = "'5.66%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1877"
$ws.Range("E11").Value = "'4.92%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08699"
$ws.Range("E12").Value = "'3.17%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03506"
$ws.Range("E13").Value = "'-1.34%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09921"
$ws.Range("E14").Value = "'-0.16%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001439"
$ws.Range("E15").Value = "'-2.56%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005715"
$ws.Range("E16").Value = "'0.89%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.459"
$ws.Range("E17").Value = "'-0.53%"

# Rows 18-27: price/volume updates
$ws.Range("D18").Value = "'2.365"
$ws.Range("E18").Value = "'6.68%"
$ws.Range("D20").Value = "'0.1346"
$ws.Range("E20").Value = "'2.05%"
$ws.Range("E21").Value = "'4.48%"
$ws.Range("E22").Value = "'1.78%"
$ws.Range("D23").Value = "'0.04607"
$ws.Range("E23").Value = "'-1.19%"
$ws.Range("E24").Value = "'15.13%"
$ws.Range("D25").Value = "'0.001230"
$ws.Range("E25").Value = "'-0.15%"
$ws.Range("D26").Value = "'0.0001401"
$ws.Range("E26").Value = "'7.68%"
$ws.Range("E27").Value = "'-42.77%"

# Rows 39-51: price/volume updates
$ws.Range("D39").Value = "'0.01828"
$ws.Range("E39").Value = "'5.04%"
$ws.Range("D40").Value = "'0.04769"
$ws.Range("E40").Value = "'1.87%"
$ws.Range("D41").Value = "'0.007495"
$ws.Range("E41").Value = "'-5.40%"
$ws.Range("D42").Value = "'0.1400"
$ws.Range("E42").Value = "'1.07%"
$ws.Range("D43").Value = "'0.007736"
$ws.Range("E43").Value = "'0.50%"
$ws.Range("D44").Value = "'0.002232"
$ws.Range("E44").Value = "'-2.64%"
$ws.Range("D45").Value = "'0.01044"
$ws.Range("E45").Value = "'7.00%"
$ws.Range("D46").Value = "'0.00006230"
$ws.Range("E46").Value = "'2.80%"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("D48").Value = "'0.0005802"
$ws.Range("E48").Value = "'0.03%"
$ws.Range("D49").Value = "'28.20"
$ws.Range("E49").Value = "'221.73%"
$ws.Range("E50").Value = "'-25.86%"
$ws.Range("D51").Value = "'0.00002101"
$ws.Range("E51").Value = "'-0.05%"
